$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at row 66. This pushes the old rows 66-94
#    (the "Bio"/"System"/"Math"/"Music"/"Probabilistic"/"Fake Algorithms" groups)
#    down to rows 67-95, making room for a 12th "Physics" entry plus a new
#    spacer row.
$ws.Rows("66:66").Insert()

# 2. The last row of the (now shifted) "Fake Algorithms" group - the Battle
#    Royale Optimization entry, now sitting at row 95 - is the row that
#    belongs up with the rest of the "Physics" group (row 65 was previously
#    just an empty "STT 12" placeholder). Move (cut) it there, carrying its
#    number formatting/highlight style along with it.
$ws.Range("C95:L95").Cut($ws.Range("C65"))

# 3. Row 66 becomes the new spacer/placeholder row for the Physics group
#    (STT counter continues to 13, no data yet).
$ws.Range("B66").Value = 13

# 4. Row 95 is now entirely empty (its data was cut out in step 2) - remove
#    it so the sheet goes back to 94 total rows, matching the rest of the
#    groups which simply shifted down by one.
$ws.Rows("95:95").Delete()

# 5. The "Fake Algorithms" group (now starting at row 86, having shifted down
#    by one row because of the insert in step 1) is renamed to "Dummy
#    Algorithms" - it lost its Battle Royale Optimization entry to the
#    Physics group above, but keeps the rest of its members.
$ws.Range("A86").Value = "Dummy Algorithms"

# 6. Restore the view/selection as it was left by the editing author.
$ws.Application.ActiveWindow.ScrollRow = 79
$sel = $ws.Range("M84")
$sel.Select()
